$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-13 with new TPM-derived values ---
# Row 2
$ws.Cells.Item(2, 7).Value = 46.105436
$ws.Cells.Item(2, 8).Value = 138.316308
$ws.Cells.Item(2, 9).Value = 0.929514990096053
$ws.Cells.Item(2, 10).Value = 0.9295149900960532
$ws.Cells.Item(2, 13).Value = 209.26237
$ws.Cells.Item(2, 14).Value = 627.78711
$ws.Cells.Item(2, 15).Value = 0.8127157202241573
$ws.Cells.Item(2, 16).Value = 0.8127157202241573
$ws.Cells.Item(2, 17).Value = 9648.132807243319
$ws.Cells.Item(2, 18).Value = 86833.19526518988
$ws.Cells.Item(2, 19).Value = 0.7554314446350642
$ws.Cells.Item(2, 20).Value = 0.7554314446350643

# Row 3
$ws.Cells.Item(3, 7).Value = 46.105436
$ws.Cells.Item(3, 8).Value = 138.316308
$ws.Cells.Item(3, 9).Value = 0.929514990096053
$ws.Cells.Item(3, 10).Value = 0.9295149900960532
$ws.Cells.Item(3, 13).Value = 0.9848756666666668
$ws.Cells.Item(3, 14).Value = 2.954627
$ws.Cells.Item(3, 15).Value = 0.003824977881910862
$ws.Cells.Item(3, 16).Value = 0.003824977881910862
$ws.Cells.Item(3, 17).Value = 45.40812201745734
$ws.Cells.Item(3, 18).Value = 408.673098157116
$ws.Cells.Item(3, 19).Value = 0.003555374278021997
$ws.Cells.Item(3, 20).Value = 0.003555374278021998

# Row 4
$ws.Cells.Item(4, 7).Value = 46.105436
$ws.Cells.Item(4, 8).Value = 138.316308
$ws.Cells.Item(4, 9).Value = 0.929514990096053
$ws.Cells.Item(4, 10).Value = 0.9295149900960532
$ws.Cells.Item(4, 13).Value = 1.763846666666667
$ws.Cells.Item(4, 14).Value = 5.291539999999999
$ws.Cells.Item(4, 15).Value = 0.006850280411451801
$ws.Cells.Item(4, 16).Value = 0.006850280411451801
$ws.Cells.Item(4, 17).Value = 81.32291960381332
$ws.Cells.Item(4, 18).Value = 731.9062764343199
$ws.Cells.Item(4, 19).Value = 0.006367438328805807
$ws.Cells.Item(4, 20).Value = 0.006367438328805808

# Row 5
$ws.Cells.Item(5, 7).Value = 46.105436
$ws.Cells.Item(5, 8).Value = 138.316308
$ws.Cells.Item(5, 9).Value = 0.929514990096053
$ws.Cells.Item(5, 10).Value = 0.9295149900960532
$ws.Cells.Item(5, 13).Value = 45.474231
$ws.Cells.Item(5, 14).Value = 136.422693
$ws.Cells.Item(5, 15).Value = 0.1766090214824801
$ws.Cells.Item(5, 16).Value = 0.1766090214824801
$ws.Cells.Item(5, 17).Value = 2096.609247019716
$ws.Cells.Item(5, 18).Value = 18869.48322317744
$ws.Cells.Item(5, 19).Value = 0.1641607328541611
$ws.Cells.Item(5, 20).Value = 0.1641607328541611

# Row 6
$ws.Cells.Item(6, 9).Value = 0.04170958390412858
$ws.Cells.Item(6, 10).Value = 0.04170958390412859
$ws.Cells.Item(6, 13).Value = 209.26237
$ws.Cells.Item(6, 14).Value = 627.78711
$ws.Cells.Item(6, 15).Value = 0.8127157202241573
$ws.Cells.Item(6, 16).Value = 0.8127157202241573
$ws.Cells.Item(6, 17).Value = 432.9350350770633
$ws.Cells.Item(6, 18).Value = 3896.415315693569
$ws.Cells.Item(6, 19).Value = 0.03389803452289378
$ws.Cells.Item(6, 20).Value = 0.03389803452289378

# Row 7
$ws.Cells.Item(7, 9).Value = 0.04170958390412858
$ws.Cells.Item(7, 10).Value = 0.04170958390412859
$ws.Cells.Item(7, 13).Value = 0.9848756666666668
$ws.Cells.Item(7, 14).Value = 2.954627
$ws.Cells.Item(7, 15).Value = 0.003824977881910862
$ws.Cells.Item(7, 16).Value = 0.003824977881910862
$ws.Cells.Item(7, 17).Value = 2.037572169783223
$ws.Cells.Item(7, 18).Value = 18.338149528049
$ws.Cells.Item(7, 19).Value = 0.0001595382358969971
$ws.Cells.Item(7, 20).Value = 0.0001595382358969971

# Row 8
$ws.Cells.Item(8, 9).Value = 0.04170958390412858
$ws.Cells.Item(8, 10).Value = 0.04170958390412859
$ws.Cells.Item(8, 13).Value = 1.763846666666667
$ws.Cells.Item(8, 14).Value = 5.291539999999999
$ws.Cells.Item(8, 15).Value = 0.006850280411451801
$ws.Cells.Item(8, 16).Value = 0.006850280411451801
$ws.Cells.Item(8, 17).Value = 3.649155930442222
$ws.Cells.Item(8, 18).Value = 32.84240337397999
$ws.Cells.Item(8, 19).Value = 0.0002857223455882573
$ws.Cells.Item(8, 20).Value = 0.0002857223455882574

# Row 9
$ws.Cells.Item(9, 9).Value = 0.04170958390412858
$ws.Cells.Item(9, 10).Value = 0.04170958390412859
$ws.Cells.Item(9, 13).Value = 45.474231
$ws.Cells.Item(9, 14).Value = 136.422693
$ws.Cells.Item(9, 15).Value = 0.1766090214824801
$ws.Cells.Item(9, 16).Value = 0.1766090214824801
$ws.Cells.Item(9, 17).Value = 94.07992365319899
$ws.Cells.Item(9, 18).Value = 846.7193128787908
$ws.Cells.Item(9, 19).Value = 0.00736628879974955
$ws.Cells.Item(9, 20).Value = 0.007366288799749551

# Row 10
$ws.Cells.Item(10, 7).Value = 1.378131333333333
$ws.Cells.Item(10, 8).Value = 4.134394
$ws.Cells.Item(10, 9).Value = 0.02778400648145685
$ws.Cells.Item(10, 10).Value = 0.02778400648145685
$ws.Cells.Item(10, 13).Value = 209.26237
$ws.Cells.Item(10, 14).Value = 627.78711
$ws.Cells.Item(10, 15).Value = 0.8127157202241573
$ws.Cells.Item(10, 16).Value = 0.8127157202241573
$ws.Cells.Item(10, 17).Value = 288.3910289845933
$ws.Cells.Item(10, 18).Value = 2595.51926086134
$ws.Cells.Item(10, 19).Value = 0.02258049883828986
$ws.Cells.Item(10, 20).Value = 0.02258049883828986

# Row 11
$ws.Cells.Item(11, 7).Value = 1.378131333333333
$ws.Cells.Item(11, 8).Value = 4.134394
$ws.Cells.Item(11, 9).Value = 0.02778400648145685
$ws.Cells.Item(11, 10).Value = 0.02778400648145685
$ws.Cells.Item(11, 13).Value = 0.9848756666666668
$ws.Cells.Item(11, 14).Value = 2.954627
$ws.Cells.Item(11, 15).Value = 0.003824977881910862
$ws.Cells.Item(11, 16).Value = 0.003824977881910862
$ws.Cells.Item(11, 17).Value = 1.357288015670889
$ws.Cells.Item(11, 18).Value = 12.215592141038
$ws.Cells.Item(11, 19).Value = 0.0001062732102624405
$ws.Cells.Item(11, 20).Value = 0.0001062732102624405

# Row 12
$ws.Cells.Item(12, 7).Value = 1.378131333333333
$ws.Cells.Item(12, 8).Value = 4.134394
$ws.Cells.Item(12, 9).Value = 0.02778400648145685
$ws.Cells.Item(12, 10).Value = 0.02778400648145685
$ws.Cells.Item(12, 13).Value = 1.763846666666667
$ws.Cells.Item(12, 14).Value = 5.291539999999999
$ws.Cells.Item(12, 15).Value = 0.006850280411451801
$ws.Cells.Item(12, 16).Value = 0.006850280411451801
$ws.Cells.Item(12, 17).Value = 2.430812358528889
$ws.Cells.Item(12, 18).Value = 21.87731122676
$ws.Cells.Item(12, 19).Value = 0.0001903282353515737
$ws.Cells.Item(12, 20).Value = 0.0001903282353515738

# Row 13
$ws.Cells.Item(13, 7).Value = 1.378131333333333
$ws.Cells.Item(13, 8).Value = 4.134394
$ws.Cells.Item(13, 9).Value = 0.02778400648145685
$ws.Cells.Item(13, 10).Value = 0.02778400648145685
$ws.Cells.Item(13, 13).Value = 45.474231
$ws.Cells.Item(13, 14).Value = 136.422693
$ws.Cells.Item(13, 15).Value = 0.1766090214824801
$ws.Cells.Item(13, 16).Value = 0.1766090214824801
$ws.Cells.Item(13, 17).Value = 62.66946260033799
$ws.Cells.Item(13, 18).Value = 564.0251634030419
$ws.Cells.Item(13, 19).Value = 0.004906906197552978
$ws.Cells.Item(13, 20).Value = 0.004906906197552979

# Row 14
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Mmrn2"
$ws.Cells.Item(14, 3).Value = "Cd93"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.049176
$ws.Cells.Item(14, 8).Value = 0.147528
$ws.Cells.Item(14, 9).Value = 0.0009914195183614249
$ws.Cells.Item(14, 10).Value = 0.000991419518361425
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 209.26237
$ws.Cells.Item(14, 14).Value = 627.78711
$ws.Cells.Item(14, 15).Value = 0.8127157202241573
$ws.Cells.Item(14, 16).Value = 0.8127157202241573
$ws.Cells.Item(14, 17).Value = 10.29068630712
$ws.Cells.Item(14, 18).Value = 92.61617676408
$ws.Cells.Item(14, 19).Value = 0.0008057422279093926
$ws.Cells.Item(14, 20).Value = 0.0008057422279093927

# Row 15
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Mmrn2"
$ws.Cells.Item(15, 3).Value = "Cd93"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.049176
$ws.Cells.Item(15, 8).Value = 0.147528
$ws.Cells.Item(15, 9).Value = 0.0009914195183614249
$ws.Cells.Item(15, 10).Value = 0.000991419518361425
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.9848756666666668
$ws.Cells.Item(15, 14).Value = 2.954627
$ws.Cells.Item(15, 15).Value = 0.003824977881910862
$ws.Cells.Item(15, 16).Value = 0.003824977881910862
$ws.Cells.Item(15, 17).Value = 0.048432245784
$ws.Cells.Item(15, 18).Value = 0.435890212056
$ws.Cells.Item(15, 19).Value = 0.00000379215772942717
$ws.Cells.Item(15, 20).Value = 0.000003792157729427171

# Row 16
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Mmrn2"
$ws.Cells.Item(16, 3).Value = "Cd93"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.049176
$ws.Cells.Item(16, 8).Value = 0.147528
$ws.Cells.Item(16, 9).Value = 0.0009914195183614249
$ws.Cells.Item(16, 10).Value = 0.000991419518361425
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 1.763846666666667
$ws.Cells.Item(16, 14).Value = 5.291539999999999
$ws.Cells.Item(16, 15).Value = 0.006850280411451801
$ws.Cells.Item(16, 16).Value = 0.006850280411451801
$ws.Cells.Item(16, 17).Value = 0.08673892368
$ws.Cells.Item(16, 18).Value = 0.7806503131199999
$ws.Cells.Item(16, 19).Value = 0.000006791501706162248
$ws.Cells.Item(16, 20).Value = 0.000006791501706162249

# Row 17
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Mmrn2"
$ws.Cells.Item(17, 3).Value = "Cd93"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.049176
$ws.Cells.Item(17, 8).Value = 0.147528
$ws.Cells.Item(17, 9).Value = 0.0009914195183614249
$ws.Cells.Item(17, 10).Value = 0.000991419518361425
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 45.474231
$ws.Cells.Item(17, 14).Value = 136.422693
$ws.Cells.Item(17, 15).Value = 0.1766090214824801
$ws.Cells.Item(17, 16).Value = 0.1766090214824801
$ws.Cells.Item(17, 17).Value = 2.236240783656
$ws.Cells.Item(17, 18).Value = 20.12616705290399
$ws.Cells.Item(17, 19).Value = 0.000175093631016443
$ws.Cells.Item(17, 20).Value = 0.000175093631016443

# --- Add new rows 14-17 for the Resolving-Mac sending cluster ---
# Row 14
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Mmrn2"
$ws.Cells.Item(14, 3).Value = "Cd93"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.049176
$ws.Cells.Item(14, 8).Value = 0.147528
$ws.Cells.Item(14, 9).Value = 0.0009914195183614249
$ws.Cells.Item(14, 10).Value = 0.000991419518361425
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 209.26237
$ws.Cells.Item(14, 14).Value = 627.78711
$ws.Cells.Item(14, 15).Value = 0.8127157202241573
$ws.Cells.Item(14, 16).Value = 0.8127157202241573
$ws.Cells.Item(14, 17).Value = 10.29068630712
$ws.Cells.Item(14, 18).Value = 92.61617676408
$ws.Cells.Item(14, 19).Value = 0.0008057422279093926
$ws.Cells.Item(14, 20).Value = 0.0008057422279093927

# Row 15
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Mmrn2"
$ws.Cells.Item(15, 3).Value = "Cd93"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.049176
$ws.Cells.Item(15, 8).Value = 0.147528
$ws.Cells.Item(15, 9).Value = 0.0009914195183614249
$ws.Cells.Item(15, 10).Value = 0.000991419518361425
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.9848756666666668
$ws.Cells.Item(15, 14).Value = 2.954627
$ws.Cells.Item(15, 15).Value = 0.003824977881910862
$ws.Cells.Item(15, 16).Value = 0.003824977881910862
$ws.Cells.Item(15, 17).Value = 0.048432245784
$ws.Cells.Item(15, 18).Value = 0.435890212056
$ws.Cells.Item(15, 19).Value = 0.00000379215772942717
$ws.Cells.Item(15, 20).Value = 0.000003792157729427171

# Row 16
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Mmrn2"
$ws.Cells.Item(16, 3).Value = "Cd93"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.049176
$ws.Cells.Item(16, 8).Value = 0.147528
$ws.Cells.Item(16, 9).Value = 0.0009914195183614249
$ws.Cells.Item(16, 10).Value = 0.000991419518361425
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 1.763846666666667
$ws.Cells.Item(16, 14).Value = 5.291539999999999
$ws.Cells.Item(16, 15).Value = 0.006850280411451801
$ws.Cells.Item(16, 16).Value = 0.006850280411451801
$ws.Cells.Item(16, 17).Value = 0.08673892368
$ws.Cells.Item(16, 18).Value = 0.7806503131199999
$ws.Cells.Item(16, 19).Value = 0.000006791501706162248
$ws.Cells.Item(16, 20).Value = 0.000006791501706162249

# Row 17
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Mmrn2"
$ws.Cells.Item(17, 3).Value = "Cd93"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.049176
$ws.Cells.Item(17, 8).Value = 0.147528
$ws.Cells.Item(17, 9).Value = 0.0009914195183614249
$ws.Cells.Item(17, 10).Value = 0.000991419518361425
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 45.474231
$ws.Cells.Item(17, 14).Value = 136.422693
$ws.Cells.Item(17, 15).Value = 0.1766090214824801
$ws.Cells.Item(17, 16).Value = 0.1766090214824801
$ws.Cells.Item(17, 17).Value = 2.236240783656
$ws.Cells.Item(17, 18).Value = 20.12616705290399
$ws.Cells.Item(17, 19).Value = 0.000175093631016443
$ws.Cells.Item(17, 20).Value = 0.000175093631016443
